{"js": "// Update the two Bibliography entries to use \"[n]\" style reference markers\n// and normalise the punctuation of each citation (commit: \"GitHub action\n// artefacts added\"):\n//   1. -> [1]   \"United Nations: ... (1948)\"      -> \"United Nations, ..., 1948.\"\n//   2. -> [2]   \"Brown, G. ed: The ...\"            -> \"G. Brown, ed., The ...\"\n//               \"... changing world. Open Book ... [New York] (2016)\"\n//               -> \"... changing world, Open Book ... [New York], 2016.\"\n\n// Locate the two paragraphs that use the \"Bibliography\" paragraph style \u2014\n// this is more robust than relying on a fixed paragraph index.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nconst bibParagraphs = paragraphs.items.filter((p) => p.style === \"Bibliography\");\nif (bibParagraphs.length < 2) {\n  throw new Error(\"Expected at least 2 Bibliography-styled paragraphs, found \" + bibParagraphs.length);\n}\nconst [entry1, entry2] = bibParagraphs;\n\n// --- Bibliography entry 1 (United Nations, 1948) ---------------------------\nfunction replaceInParagraph(paragraph, searchText, replacement) {\n  const found = paragraph.search(searchText, { matchCase: true });\n  found.load(\"items\");\n  return found;\n}\n\nlet found1 = replaceInParagraph(entry1, \"1.\", \"[1]\");\nlet found2 = replaceInParagraph(\n  entry1,\n  \"United Nations: Universal Declaration of Human Rights. (1948)\",\n  \"United Nations, Universal Declaration of Human Rights, 1948.\"\n);\nawait context.sync();\nif (found1.items.length !== 1) {\n  throw new Error('Expected exactly one match for \"1.\" in entry 1, found ' + found1.items.length);\n}\nif (found2.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the 1948 citation text, found \" + found2.items.length);\n}\nfound1.items[0].insertText(\"[1]\", Word.InsertLocation.replace);\nfound2.items[0].insertText(\n  \"United Nations, Universal Declaration of Human Rights, 1948.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Bibliography entry 2 (Brown, 2016) -------------------------------------\nlet found3 = replaceInParagraph(entry2, \"2.\", \"[2]\");\nlet found4 = replaceInParagraph(entry2, \"Brown, G. ed: The\", \"G. Brown, ed., The\");\nlet found5 = replaceInParagraph(\n  entry2,\n  \"in the 21st century, a living document in a changing world. Open Book Publishers ; NYU Global Institute for Advanced Study, Cambridge, [New York] (2016)\",\n  \"in the 21st century, a living document in a changing world, Open Book Publishers ; NYU Global Institute for Advanced Study, Cambridge, [New York], 2016.\"\n);\nawait context.sync();\nif (found3.items.length !== 1) {\n  throw new Error('Expected exactly one match for \"2.\" in entry 2, found ' + found3.items.length);\n}\nif (found4.items.length !== 1) {\n  throw new Error('Expected exactly one match for \"Brown, G. ed: The\", found ' + found4.items.length);\n}\nif (found5.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the 2016 citation tail text, found \" + found5.items.length);\n}\nfound3.items[0].insertText(\"[2]\", Word.InsertLocation.replace);\nfound4.items[0].insertText(\"G. Brown, ed., The\", Word.InsertLocation.replace);\nfound5.items[0].insertText(\n  \"in the 21st century, a living document in a changing world, Open Book Publishers ; NYU Global Institute for Advanced Study, Cambridge, [New York], 2016.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Update the two Bibliography entries to use \"[n]\" style reference markers\n# and normalise the punctuation of each citation (commit: \"GitHub action\n# artefacts added\"):\n#   1. -> [1]   \"United Nations: ... (1948)\"   -> \"United Nations, ..., 1948.\"\n#   2. -> [2]   \"Brown, G. ed: The ...\"         -> \"G. Brown, ed., The ...\"\n#               \"... changing world. Open Book ... [New York] (2016)\"\n#               -> \"... changing world, Open Book ... [New York], 2016.\"\n\n$d = $word.ActiveDocument\n\n# Locate the two paragraphs that use the \"Bibliography\" paragraph style --\n# more robust than relying on a fixed paragraph index.\n$bibParaIndexes = New-Object System.Collections.ArrayList\n$paraCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $paraCount; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Style.NameLocal -eq \"Bibliography\") {\n        [void]$bibParaIndexes.Add($i)\n    }\n}\nif ($bibParaIndexes.Count -lt 2) {\n    throw \"Expected at least 2 Bibliography-styled paragraphs, found $($bibParaIndexes.Count)\"\n}\n$entry1Index = $bibParaIndexes[0]\n$entry2Index = $bibParaIndexes[1]\n\nfunction Replace-TextInParagraph($paraIndex, $searchText, $replaceText) {\n    $rng = $word.ActiveDocument.Paragraphs.Item($paraIndex).Range\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $found = $find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n    if (-not $found) {\n        throw \"Could not find '$searchText' in paragraph $paraIndex\"\n    }\n}\n\n# --- Bibliography entry 1 (United Nations, 1948) ---------------------------\nReplace-TextInParagraph $entry1Index \"1.\" \"[1]\"\nReplace-TextInParagraph $entry1Index \"United Nations: Universal Declaration of Human Rights. (1948)\" \"United Nations, Universal Declaration of Human Rights, 1948.\"\n\n# --- Bibliography entry 2 (Brown, 2016) -------------------------------------\nReplace-TextInParagraph $entry2Index \"2.\" \"[2]\"\nReplace-TextInParagraph $entry2Index \"Brown, G. ed: The\" \"G. Brown, ed., The\"\nReplace-TextInParagraph $entry2Index \"in the 21st century, a living document in a changing world. Open Book Publishers ; NYU Global Institute for Advanced Study, Cambridge, [New York] (2016)\" \"in the 21st century, a living document in a changing world, Open Book Publishers ; NYU Global Institute for Advanced Study, Cambridge, [New York], 2016.\"\n"}
